$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.387524366378784
$ws.Range("B1").Value = 2.577778100967407
$ws.Range("C1").Value = 6.601644992828369
$ws.Range("D1").Value = 2.407756328582764
$ws.Range("E1").Value = 1.208657026290894
